$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A17").Value = "Cls Coleslaw 2 x 1Kg"
$ws.Range("A18").Value = "Cls Coleslaw 6 x 400g"
$ws.Range("A19").Value = "Cls Coleslaw 6 x 800g"
$ws.Range("A20").Value = "Cls Coleslaw 8 x 250g"
$ws.Range("A29").Value = "Cls Greek Salad 2 x 1Kg"
$ws.Range("A36").Value = "Cls Pasta Sld 2 x 1Kg"
$ws.Range("A37").Value = "Cls Pasta Sld 6 x 400g"
$ws.Range("A38").Value = "Cls Pasta Sld 6 x 800g"
$ws.Range("A39").Value = "Cls Pasta Sld 6x400g(TAS)"
$ws.Range("A40").Value = "Cls Pasta Sld 6x800g(TAS)"
$ws.Range("A41").Value = "Cls Pasta Sld 8 x 250g"
$ws.Range("A42").Value = "Cls Pasta Sld 8x250g(TAS)"
$ws.Range("A46").Value = "Cls Potato Sld 2 x 1Kg"
$ws.Range("A47").Value = "Cls Potato Sld 6 x 400g"
$ws.Range("A48").Value = "Cls Potato Sld 6 x 800g"
$ws.Range("A49").Value = "Cls Potato Sld 8 x 250g"
$ws.Range("A60").Value = "Cls YLC SeafdSld 3x1kg"
$ws.Range("A128").Value = "WW Coleslaw 12 x 110g"
$ws.Range("A129").Value = "WW Coleslaw 6 x 250g"
$ws.Range("A130").Value = "WW Coleslaw 6 x 400g"
$ws.Range("A131").Value = "WW Coleslaw 6 x 800g"
$ws.Range("A142").Value = "WW Pasta Sld 12 x 110g"
$ws.Range("A143").Value = "WW Pasta Sld 6 x 250g"
$ws.Range("A144").Value = "WW Pasta Sld 6 x 400g"
$ws.Range("A145").Value = "WW Pasta Sld 6 x 800g"
$ws.Range("A148").Value = "WW Potato Sld 6 x 250g"
$ws.Range("A149").Value = "WW Potato Sld 6 x 400g"
$ws.Range("A150").Value = "WW Potato Sld 6 x 800g"
$ws.Range("A151").Value = "WW Psta Sld 12x110g (TAS)"
$ws.Range("A152").Value = "WW Psta Sld 6x250g (TAS)"
$ws.Range("A153").Value = "WW Psta Sld 6x400g (TAS)"
$ws.Range("A154").Value = "WW Psta Sld 6x800g (TAS)"

$ws.Range("C12").Select()
